$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.16936194896698
$ws.Range("B1").Value = 2.310781717300415
$ws.Range("C1").Value = 3.16610312461853
$ws.Range("D1").Value = 1.435409665107727
$ws.Range("E1").Value = 1.109761238098145
